$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append rows 35-42 ---
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "14:21:34", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:21:37", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:21:55", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:22:10", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:22:22", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:22:24", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "14:22:44", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:22:46", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 35
$endRow = $startRow + $proximityRows.Count - 1

# Force column A (date strings like "2026-02-01") to be stored as literal
# text instead of being auto-converted to a date serial number, then drop
# the temporary formatting again so the new cells end up unstyled, same
# as all the pre-existing rows on this sheet.
$proximityDateRangeAddr = "A" + $startRow + ":A" + $endRow
$proximityDateRange = $proximity.Range($proximityDateRangeAddr)
$proximityDateRange.NumberFormat = "@"

for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $proximityRows[$i]
    $proximity.Cells.Item($r, 1).Value = $rowData[0]
    $proximity.Cells.Item($r, 2).Value = $rowData[1]
    $proximity.Cells.Item($r, 3).Value = $rowData[2]
    $proximity.Cells.Item($r, 4).Value = $rowData[3]
    $proximity.Cells.Item($r, 5).Value = $rowData[4]
    $proximity.Cells.Item($r, 6).Value = $rowData[5]
}

$proximityDateRange.ClearFormats()

# --- Camera sheet: append rows 21-25 ---
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "14:21:36", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:21:56", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:22:09", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:22:24", "14:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "14:22:45", "14:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow = 21
$endRow = $startRow + $cameraRows.Count - 1

$cameraDateRangeAddr = "A" + $startRow + ":A" + $endRow
$cameraDateRange = $camera.Range($cameraDateRangeAddr)
$cameraDateRange.NumberFormat = "@"

for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $cameraRows[$i]
    $camera.Cells.Item($r, 1).Value = $rowData[0]
    $camera.Cells.Item($r, 2).Value = $rowData[1]
    $camera.Cells.Item($r, 3).Value = $rowData[2]
    $camera.Cells.Item($r, 4).Value = $rowData[3]
    $camera.Cells.Item($r, 5).Value = $rowData[4]
    $camera.Cells.Item($r, 6).Value = $rowData[5]
}

$cameraDateRange.ClearFormats()
